# Update countries & provincias Spain
# - Reorder "Paises Bajos" ahead of "Marruecos"/"Catar" (rows 33-35), with
#   refreshed case/death counters for the three countries.
# - Refresh the "Datos actualizados" timestamp string.
# - Refresh numeric COVID counters for several other country rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder Paises Bajos / Marruecos / Catar (rows 33-35) + updated figures ---
# Row 33: now Paises Bajos (updated counters)
$ws.Range("A33").Value = "Paises Bajos"
$ws.Range("B33").Value = 127922
$ws.Range("C33").Value = 3825
$ws.Range("D33").Value = 0
$ws.Range("E33").Value = 0
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 9
$ws.Range("H33").Value = 6428

# Row 34: now Marruecos (figures carried over from former row 33)
$ws.Range("A34").Value = "Marruecos"
$ws.Range("B34").Value = 126044
$ws.Range("C34").Value = 0
$ws.Range("D34").Value = 104136
$ws.Range("E34").Value = 19679
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 0
$ws.Range("H34").Value = 2229

# Row 35: now Catar (figures carried over from former row 34)
$ws.Range("A35").Value = "Catar"
$ws.Range("B35").Value = 125959
$ws.Range("C35").Value = 0
$ws.Range("D35").Value = 122911
$ws.Range("E35").Value = 2834
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 0
$ws.Range("H35").Value = 214

# --- Estados Unidos (row 4) ---
$ws.Range("B4").Value = 7499991
$ws.Range("C4").Value = 5320
$ws.Range("D4").Value = 4737374
$ws.Range("E4").Value = 2549891
$ws.Range("G4").Value = 66
$ws.Range("H4").Value = 212726

# --- India (row 5) ---
$ws.Range("B5").Value = 6399329
$ws.Range("C5").Value = 7369
$ws.Range("D5").Value = 5353120
$ws.Range("E5").Value = 946372
$ws.Range("G5").Value = 33
$ws.Range("H5").Value = 99837

# --- Alemania (row 20) ---
$ws.Range("B20").Value = 335578
$ws.Range("C20").Value = 481
$ws.Range("D20").Value = 320348
$ws.Range("E20").Value = 10407
$ws.Range("G20").Value = 29
$ws.Range("H20").Value = 4823

# --- Ecuador (row 25) ---
$ws.Range("B25").Value = 296364
$ws.Range("C25").Value = 834
$ws.Range("E25").Value = 27276
$ws.Range("G25").Value = 2
$ws.Range("H25").Value = 9588

# --- Guatemala (row 50) ---
$ws.Range("B50").Value = 79421
$ws.Range("C50").Value = 402
$ws.Range("D50").Value = 74982
$ws.Range("E50").Value = 3595
$ws.Range("G50").Value = 5
$ws.Range("H50").Value = 844

# --- row 105 ---
$ws.Range("E105").Value = 1799
$ws.Range("G105").Value = 1
$ws.Range("H105").Value = 345

# --- Guinea (row 107) ---
$ws.Range("B107").Value = 9852
$ws.Range("C107").Value = 41
$ws.Range("D107").Value = 8661
$ws.Range("E107").Value = 1114

# --- row 163 ---
$ws.Range("B163").Value = 1680
$ws.Range("C163").Value = 41
$ws.Range("D163").Value = 926
$ws.Range("E163").Value = 715
$ws.Range("G163").Value = 1
$ws.Range("H163").Value = 39

# --- row 168 ---
$ws.Range("D168").Value = 1020
$ws.Range("E168").Value = 40

# --- row 179 ---
$ws.Range("D179").Value = 432
$ws.Range("E179").Value = 40

# --- Updated timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 2 de Octubre de 2020 a las 15:35"
